$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New entries in rows 12-14 (dates, hours, activity text)
# Copy the existing date formatting from A11 so no new number format/style gets created
$ws.Range("A11").Copy()
$ws.Range("A12:A14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A12").Value = 44174
$ws.Range("B12").Value = 5
$ws.Range("D12").Value = "Dokus und Tutorials über Vue.js gelesen/gesehen"

$ws.Range("A13").Value = 44182
$ws.Range("B13").Value = 6
$ws.Range("D13").Value = "Dokus und Tutorials über Vue.js gelesen/gesehen"

$ws.Range("A14").Value = 44183
$ws.Range("B14").Value = 6
$ws.Range("D14").Value = "Meeting + WebUI"

# Update current selection to match author's last edit location
$ws.Range("D17").Select()

$wb.Save()
